# Insert a new weekly price record for "Berenjena" (Vega Central Mapocho de
# Santiago) as row 314, pushing the previously existing rows 314-327 down to
# 315-328 (the sheet's dimension grows from A1:R327 to A1:R328).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data down by inserting a blank row at 314.
$ws.Rows.Item(314).Insert()

# Populate the newly inserted row with this week's data.
$ws.Range("A314").Value = 9
$ws.Range("B314").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C314").Value = "Metropolitana"
$ws.Range("D314").Value = 44939
$ws.Range("E314").Value = 13
$ws.Range("F314").Value = 100112001
$ws.Range("G314").Value = "Berenjena"
$ws.Range("H314").Value = "Sin especificar"
$ws.Range("I314").Value = "Primera"
$ws.Range("J314").Value = 70
$ws.Range("K314").Value = 14000
$ws.Range("L314").Value = 15000
$ws.Range("M314").Value = 14500
$ws.Range("N314").Value = "`$/caja 60 unidades"
$ws.Range("O314").Value = "Región de O'Higgins"
$ws.Range("P314").Value = 242
$ws.Range("Q314").Value = 60
$ws.Range("R314").Value = "Hortaliza"
